$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 24.60000000000041
$ws.Range("F3").Value = 24.60000000000041
$ws.Range("F4").Value = 24.60000000000041
$ws.Range("F5").Value = 24.60000000000041
$ws.Range("F6").Value = 24.60000000000041
$ws.Range("F7").Value = 24.99000000000047
$ws.Range("F8").Value = 24.99000000000047
$ws.Range("F9").Value = 24.99000000000047
$ws.Range("F10").Value = 24.99000000000047
$ws.Range("F11").Value = 24.99000000000047
$ws.Range("F12").Value = 24.99000000000047
$ws.Range("H2").Value = 0.09300270366551344
$ws.Range("H3").Value = 0.2460937912464798
$ws.Range("H4").Value = 0.6125823743994973
$ws.Range("H5").Value = 0.1976434392364035
$ws.Range("H6").Value = 0.5612905178853518
$ws.Range("H7").Value = 0.4217887497228616
$ws.Range("H8").Value = 0.1007531249075665
$ws.Range("H9").Value = 0.1831387382462031
$ws.Range("H10").Value = 0.1721052843877694
$ws.Range("H11").Value = 0.4755067224826074
$ws.Range("H12").Value = 0.05289958558096775
$ws.Range("I2").Value = 0.09300270366551344
$ws.Range("I3").Value = 0.2460937912464798
$ws.Range("I4").Value = 0.6125823743994973
$ws.Range("I5").Value = 0.1976434392364035
$ws.Range("I6").Value = 0.5612905178853518
$ws.Range("I7").Value = 0.4217887497228616
$ws.Range("I8").Value = 0.1007531249075665
$ws.Range("I9").Value = 0.1831387382462031
$ws.Range("I10").Value = 0.1721052843877694
$ws.Range("I11").Value = 0.4755067224826074
$ws.Range("I12").Value = 0.05289958558096775
$ws.Range("L2").Value = 7.448497245427472
$ws.Range("L3").Value = 5.226773335802929
$ws.Range("L4").Value = 3.433455085108109
$ws.Range("L5").Value = 5.663408465039712
$ws.Range("L6").Value = 3.490477283058333
$ws.Range("L7").Value = 4.775078028849081
$ws.Range("L8").Value = 6.272130088012329
$ws.Range("L9").Value = 7.698793850822242
$ws.Range("L10").Value = 5.840279630001479
$ws.Range("L11").Value = 4.253473096690889
$ws.Range("L12").Value = 6.894700057883911
$ws.Range("M2").Value = "[-1.0669984548577585, 15.963992945712702]"
$ws.Range("M3").Value = "[-2.8616155143492454, 13.315162185955103]"
$ws.Range("M4").Value = "[-6.014939154572211, 12.88184932478843]"
$ws.Range("M5").Value = "[-2.4188577580121784, 13.745674688091603]"
$ws.Range("M6").Value = "[-5.319975849817037, 12.300930415933703]"
$ws.Range("M7").Value = "[-4.813944218318908, 14.36410027601707]"
$ws.Range("M8").Value = "[-1.0702675576428238, 13.614527733667481]"
$ws.Range("M9").Value = "[-3.471161853327991, 18.868749554972474]"
$ws.Range("M10").Value = "[-2.7477835518423444, 14.428342811845303]"
$ws.Range("M11").Value = "[-4.988169531869912, 13.49511572525169]"
$ws.Range("M12").Value = "[-0.5963840422634572, 14.38578415803128]"
$ws.Range("N2").Value = 0.08490685239776119
$ws.Range("N3").Value = 0.1997019285030197
$ws.Range("N4").Value = 0.4680209109004669
$ws.Range("N5").Value = 0.1650262692474893
$ws.Range("N6").Value = 0.4290981664276425
$ws.Range("N7").Value = 0.3212383210542038
$ws.Range("N8").Value = 0.09220906992312949
$ws.Range("N9").Value = 0.1719097954229016
$ws.Range("N10").Value = 0.1775821631804584
$ws.Range("N11").Value = 0.3588757116227452
$ws.Range("N12").Value = 0.07033622441743037
$ws.Range("O2").Value = 0.08490685239776119
$ws.Range("O3").Value = 0.1997019285030197
$ws.Range("O4").Value = 0.4680209109004669
$ws.Range("O5").Value = 0.1650262692474893
$ws.Range("O6").Value = 0.4290981664276425
$ws.Range("O7").Value = 0.3212383210542038
$ws.Range("O8").Value = 0.09220906992312949
$ws.Range("O9").Value = 0.1719097954229016
$ws.Range("O10").Value = 0.1775821631804584
$ws.Range("O11").Value = 0.3588757116227452
$ws.Range("O12").Value = 0.07033622441743037
$ws.Range("P2").Value = -1.446579199851156
$ws.Range("P3").Value = -1.861684535460618
$ws.Range("P4").Value = -2.578684660604234
$ws.Range("P5").Value = -1.861684535460618
$ws.Range("P6").Value = -1.773631888513156
$ws.Range("P7").Value = 1.742184514603348
$ws.Range("P8").Value = 1.767342413731195
$ws.Range("P9").Value = 2.283079345852042
$ws.Range("P10").Value = 2.333395144107734
$ws.Range("P11").Value = 1.213868632918579
$ws.Range("P12").Value = 0.761026448617347
$ws.Range("Q2").Value = "[-2.9937899962136965, 0.10063159651138465]"
$ws.Range("Q3").Value = "[-4.937237703839816, 1.2138686329185795]"
$ws.Range("Q4").Value = "[-5.7171325768030465, 0.5597632555945777]"
$ws.Range("Q5").Value = "[-4.9435271786217765, 1.2201581077005397]"
$ws.Range("Q6").Value = "[-4.9120798047119685, 1.3648160276856558]"
$ws.Range("Q7").Value = "[-1.3836844520315408, 4.8680534812382374]"
$ws.Range("Q8").Value = "[0.10692107129334616, 3.4277637561690444]"
$ws.Range("Q9").Value = "[-0.8365001460008852, 5.402658837704969]"
$ws.Range("Q10").Value = "[-0.8050527720910781, 5.471843060306546]"
$ws.Range("Q11").Value = "[-1.9182898084982725, 4.34602707433543]"
$ws.Range("Q12").Value = "[-1.3899739268135018, 2.912026824048196]"
$ws.Range("R2").Value = 0.06615949986271641
$ws.Range("R3").Value = 0.229131560696914
$ws.Range("R4").Value = 0.1049091022515263
$ws.Range("R5").Value = 0.2300679234632435
$ws.Range("R6").Value = 0.2610479266730927
$ws.Range("R7").Value = 0.2675830802733412
$ws.Range("R8").Value = 0.03748855103701532
$ws.Range("R9").Value = 0.1474349381985349
$ws.Range("R10").Value = 0.1412557942207542
$ws.Range("R11").Value = 0.4391445814593646
$ws.Range("R12").Value = 0.4797769782272965
$ws.Range("S2").Value = 0.06615949986271641
$ws.Range("S3").Value = 0.229131560696914
$ws.Range("S4").Value = 0.1049091022515263
$ws.Range("S5").Value = 0.2300679234632435
$ws.Range("S6").Value = 0.2610479266730927
$ws.Range("S7").Value = 0.2675830802733412
$ws.Range("S8").Value = 0.03748855103701532
$ws.Range("S9").Value = 0.1474349381985349
$ws.Range("S10").Value = 0.1412557942207542
$ws.Range("S11").Value = 0.4391445814593646
$ws.Range("S12").Value = 0.4797769782272965
$ws.Range("T2").Value = 14.24730820888201
$ws.Range("T3").Value = 11.04934999506458
$ws.Range("T4").Value = 14.73444705706061
$ws.Range("T5").Value = 13.19115768773635
$ws.Range("T6").Value = 14.78340734097695
$ws.Range("T7").Value = 14.73201737330555
$ws.Range("T8").Value = 12.33150325374959
$ws.Range("T9").Value = 14.94667499551281
$ws.Range("T10").Value = 10.53533199829142
$ws.Range("T11").Value = 14.38267596542229
$ws.Range("T12").Value = 13.36339994101833
$ws.Range("U2").Value = "[9.45015758753572, 19.044458830228308]"
$ws.Range("U3").Value = "[6.637670947725473, 15.461029042403682]"
$ws.Range("U4").Value = "[9.862808973654516, 19.606085140466707]"
$ws.Range("U5").Value = "[8.756144913523485, 17.626170461949222]"
$ws.Range("U6").Value = "[10.140874538937773, 19.425940143016135]"
$ws.Range("U7").Value = "[9.496531835298807, 19.967502911312298]"
$ws.Range("U8").Value = "[8.18429466140263, 16.478711846096544]"
$ws.Range("U9").Value = "[9.042320042447486, 20.85102994857813]"
$ws.Range("U10").Value = "[6.150721261517919, 14.919942735064911]"
$ws.Range("U11").Value = "[9.403697302808496, 19.361654628036092]"
$ws.Range("U12").Value = "[9.48110124068032, 17.24569864135633]"
$ws.Range("V2").Value = 0.0000003322869319699606
$ws.Range("V3").Value = 0.000007932155460022017
$ws.Range("V4").Value = 0.0000002281705040463322
$ws.Range("V5").Value = 0.0000003224403999357861
$ws.Range("V6").Value = 0.00000007579514704758594
$ws.Range("V7").Value = 0.0000009705602839193972
$ws.Range("V8").Value = 0.0000003243900614791073
$ws.Range("V9").Value = 0.000006619246198713569
$ws.Range("V10").Value = 0.00001566785932105574
$ws.Range("V11").Value = 0.000000581034026003735
$ws.Range("V12").Value = 0.00000001281814010489768
$ws.Range("W2").Value = 0.0000003322869319699606
$ws.Range("W3").Value = 0.000007932155460022017
$ws.Range("W4").Value = 0.0000002281705040463322
$ws.Range("W5").Value = 0.0000003224403999357861
$ws.Range("W6").Value = 0.00000007579514704758594
$ws.Range("W7").Value = 0.0000009705602839193972
$ws.Range("W8").Value = 0.0000003243900614791073
$ws.Range("W9").Value = 0.000006619246198713569
$ws.Range("W10").Value = 0.00001566785932105574
$ws.Range("W11").Value = 0.000000581034026003735
$ws.Range("W12").Value = 0.00000001281814010489768
$ws.Range("X2").Value = 5.663663663663758
$ws.Range("X3").Value = 7.288888888889012
$ws.Range("X4").Value = 10.09609609609626
$ws.Range("X5").Value = 7.288888888889012
$ws.Range("X6").Value = 6.944144144144257
$ws.Range("X7").Value = 18.06084084084118
$ws.Range("X8").Value = 17.96078078078111
$ws.Range("X9").Value = 15.90954954954985
$ws.Range("X10").Value = 15.70942942942972
$ws.Range("X11").Value = 20.16210210210248
$ws.Range("X12").Value = 21.9631831831836
$ws.Range("Y2").Value = -0.3939939939939983
$ws.Range("Y3").Value = -4.752552552552633
$ws.Range("Y4").Value = -2.19159159159163
$ws.Range("Y5").Value = -4.777177177177251
$ws.Range("Y6").Value = -5.343543543543634
$ws.Range("Y7").Value = 5.628378378378487
$ws.Range("Y8").Value = 11.35681681681702
$ws.Range("Y9").Value = 3.502102102102169
$ws.Range("Y10").Value = 3.226936936936999
$ws.Range("Y11").Value = 7.704624624624769
$ws.Range("Y12").Value = 13.4080480480483
$ws.Range("Z2").Value = 11.72132132132151
$ws.Range("Z3").Value = 19.33033033033066
$ws.Range("Z4").Value = 22.38378378378415
$ws.Range("Z5").Value = 19.35495495495528
$ws.Range("Z6").Value = 19.23183183183215
$ws.Range("Z7").Value = 30.49330330330388
$ws.Range("Z8").Value = 24.56474474474521
$ws.Range("Z9").Value = 28.31699699699752
$ws.Range("Z10").Value = 28.19192192192245
$ws.Range("Z11").Value = 32.61957957958019
$ws.Range("Z12").Value = 30.51831831831889
